# Updated Project Summary _Intermediate_Interface_V4
#
# The document ends with a trailing, otherwise-empty paragraph that only
# carries the `_GoBack` bookmark. This script expands that single
# paragraph into four paragraphs:
#   1) "//"
#   2) the "story" paragraph (the _GoBack bookmark ends up in the middle
#      of it, right where it originally sat)
#   3) "So I envisioned ..."
#   4) "But than it evolved ..."
#
# NOTE: Paragraph.Index is not re-evaluated after a mutation on an
# already-fetched Paragraph object in this host, so every paragraph is
# re-fetched by its (plain-integer) position straight from
# $d.Paragraphs.Item(...) right after each structural edit instead of
# being cached.

$d = $word.ActiveDocument

# The bookmark-only paragraph is always the last paragraph in the body.
$bmIndex = $d.Paragraphs.Count
$bmPara = $d.Paragraphs.Item($bmIndex)

# Carve out three more empty paragraphs from it so we end up with four
# paragraphs total, in document order:
#   slashIndex  -> "//"
#   storyIndex  -> "story" paragraph (keeps the _GoBack bookmark)
#   storyIndex+1-> "So I envisioned ..."
#   storyIndex+2-> "But than it evolved ..."

# InsertParagraphBefore() adds a brand new empty paragraph just *before*
# $bmPara; the bookmark paragraph itself (and its bookmark) shifts one
# slot later, to $bmIndex + 1.
$bmPara.Range.InsertParagraphBefore()
$slashIndex = $bmIndex
$storyIndex = $bmIndex + 1

# InsertParagraphAfter() adds a new empty paragraph right after the
# paragraph whose Range it's called on, without disturbing that
# paragraph's own contents (the bookmark stays put in $storyIndex).
$storyPara = $d.Paragraphs.Item($storyIndex)
$storyPara.Range.InsertParagraphAfter()
$storyPara = $d.Paragraphs.Item($storyIndex)
$storyPara.Range.InsertParagraphAfter()

$envisionedIndex = $storyIndex + 1
$evolvedIndex = $storyIndex + 2

# --- Paragraph 1: "//" ------------------------------------------------
$slashPara = $d.Paragraphs.Item($slashIndex)
$slashPara.Range.Text = "//"

# --- Paragraph 2: the "story" paragraph --------------------------------
# Build the sentence as one run per sentence-fragment, in reading order,
# growing the insertion point forward each time so every InsertAfter
# lands on a fresh zero-length Range (keeps the runs distinct instead of
# merging them into one run).
$storyPara = $d.Paragraphs.Item($storyIndex)
$pos = $storyPara.Range.Start

$fragments = @(
    "The story of my project begins when I realized how crazy ",
    "inconvenient",
    " it was to ",
    "have ",
    "more ",
    "than",
    " 2 ",
    "application",
    " on your desktop at one time."
)

foreach ($fragment in $fragments) {
    $ins = $d.Range($pos, $pos)
    $ins.InsertAfter($fragment)
    $pos = $ins.End
}

# Re-seat the _GoBack bookmark between "have " and "more ", i.e. right
# after the four leading fragments. Bookmarks.Add with an existing name
# relocates that bookmark rather than creating a duplicate.
$prefixLength = ("The story of my project begins when I realized how crazy " + `
                  "inconvenient" + `
                  " it was to " + `
                  "have ").Length
$storyPara = $d.Paragraphs.Item($storyIndex)
$bmPos = $storyPara.Range.Start + $prefixLength
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Paragraph 3: "So I envisioned ..." --------------------------------
$envisionedPara = $d.Paragraphs.Item($envisionedIndex)
$envisionedPara.Range.Text = "So I envisioned an application where you are able to draw and drag your applications to optimize your real state. "

# --- Paragraph 4: "But than it evolved ..." -----------------------------
$evolvedPara = $d.Paragraphs.Item($evolvedIndex)
$evolvedPara.Range.Text = "But than it evolved, I was wondering that if I can control other applications on the desktop will I be able to make a quick executable that you click twice, and it opens, formats and sets up your working desktop."
